# Applies the "Kurseinheit eingetragen Handbuch fertigestellt" edit:
#  - Splits several table-cell paragraphs into extra runs so that the
#    English/technical words (tkinter, pygame, Perzeptron, Machine) sit in
#    their own <w:r>, matching the proofed/spell-checked run layout.
#  - Fills in the missing "18" Kurseinheit number for the pygame row.
#  - Bumps the cached TIME field in the header from 09.03.2024 to 16.03.2024.

$d = $word.ActiveDocument

function Split-RunAt($startOffset, $length) {
    # Forces Word to break the run containing [startOffset, startOffset+length)
    # into (before | target | after) runs without changing visible formatting:
    # toggling Bold on then back off is enough to force a run split, and the
    # final off-state leaves the run properties identical to their neighbours.
    $sub = $d.Range($startOffset, $startOffset + $length)
    $sub.Bold = 1
    $sub.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "GUI mit tkinter – Computer errät Zahl" -> isolate "tkinter"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("GUI mit tkinter", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$wordStart = $rng.Start + "GUI mit ".Length
Split-RunAt $wordStart "tkinter".Length

# ---------------------------------------------------------------------
# 2) "Bibliothek pygame – Traffic Game" -> isolate "pygame"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Bibliothek pygame", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$wordStart = $rng.Start + "Bibliothek ".Length
Split-RunAt $wordStart "pygame".Length

# ---------------------------------------------------------------------
# 3) Kurseinheit "18" for the pygame/Traffic Game row (empty cell)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Bibliothek pygame", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$cell = $rng.Cells.Item(1)
$rowIdx = $cell.RowIndex
$table = $cell.Tables.Item(1)
$numCell = $table.Cell($rowIdx, 3)
$numCell.Range.InsertAfter("18")

# ---------------------------------------------------------------------
# 4) "Künstliche Intelligenz (KI) – Perzeptron " -> isolate "Perzeptron"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Künstliche Intelligenz (KI) – Perzeptron", $true, $false, `
                   $false, $false, $false, $true, 1, $false, "", 0)
$wordStart = $rng.Start + "Künstliche Intelligenz (KI) – ".Length
Split-RunAt $wordStart "Perzeptron".Length

# ---------------------------------------------------------------------
# 5) "KI – Machine Learning Modelle trainieren und nutzen" -> isolate "Machine"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Machine Learning Modelle trainieren und nutzen", $true, `
                   $false, $false, $false, $false, $true, 1, $false, "", 0)
Split-RunAt $rng.Start "Machine".Length

# ---------------------------------------------------------------------
# 6) Header date field cached text: 09.03.2024 -> 16.03.2024
# ---------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute("09.03.2024", $true, $false, $false, $false, `
                             $false, $true, 1, $false, "16.03.2024", 2)
}

Write-Output "done"
